# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker "JUDITH MARGARITA DIAZ AGAMEZ" (CC 45528478) and her six
# "Periodo Mora" rows (2507..2502) are removed from the account-statement
# table. Deleting those rows shifts the remaining two worker rows
# (JOHANA PAOLA POLANCO PIANETA and ISABELA MARIA PADILLA PEREZ) up into
# rows 16-17, and the signature block up into rows 22-23.
#
# The summary figures (total mora value, worker count, period count) are
# then updated to reflect the smaller remaining data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the deleted worker's six data rows (rows 16-21); everything below
# shifts up automatically (rows, merged cells, shared strings).
$ws.Rows("16:21").Delete()

# Update the summary block to match the reduced data set.
$ws.Range("E11").Value = 56774
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2
